# Y4_B2526_Excuses.xlsx -- attendance app sync
# Updates the Student ID (col A) on the existing excuse log rows (new IDs
# issued by the system) and appends the newly-logged excuses as rows 25-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated Student IDs for existing rows 2-24 ------------------------
$idUpdates = @{
    2  = "212314"
    3  = "212315"
    4  = "212320"
    5  = "212324"
    6  = "212333"
    7  = "212345"
    8  = "212346"
    9  = "212347"
    10 = "212362"
    11 = "212368"
    12 = "212370"
    13 = "212373"
    14 = "212377"
    15 = "212414"
    16 = "212431"
    17 = "212433"
    18 = "212440"
    19 = "212445"
    20 = "212453"
    21 = "212455"
    22 = "212477"
    23 = "212479"
    24 = "212501"
}

# --- new rows 25-35 (same subject/date/time/type/user for all) ---------
$newRows = @{
    25 = "212511"
    26 = "212587"
    27 = "212588"
    28 = "212592"
    29 = "212598"
    30 = "221758"
    31 = "221884"
    32 = "223003"
    33 = "223004"
    34 = "223005"
    35 = "223006"
}

$subject = "general surgery"
$logDate = "09/09/2025"
$logTime = 0.4375
$logType = "Excuse"
$logUser = "System"

# Keep a clean "text formatted" cell on hand: writing a numeric-looking
# string into a General-formatted cell makes Excel coerce it to a number,
# so after each write we re-stamp the cell's original number format /
# style by pasting formats from a still-untouched cell of the same row
# parity (style 2/3 on even rows, 4/5 on odd rows).
$ws.Range("B2").Copy() | Out-Null
$evenFormat = $ws.Range("B2")
$ws.Range("B3").Copy() | Out-Null
$oddFormat = $ws.Range("B3")

foreach ($row in ($idUpdates.Keys | Sort-Object)) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $idUpdates[$row]

    if ($row % 2 -eq 0) {
        $evenFormat.Copy() | Out-Null
    } else {
        $oddFormat.Copy() | Out-Null
    }
    $cell.PasteSpecial(-4122) | Out-Null
}

foreach ($row in ($newRows.Keys | Sort-Object)) {
    # Clone the whole row's look from the matching-parity template row
    # (row 2 = even style, row 3 = odd style) before writing any values.
    if ($row % 2 -eq 0) {
        $ws.Range("A2:F2").Copy() | Out-Null
    } else {
        $ws.Range("A3:F3").Copy() | Out-Null
    }
    $destRow = $ws.Range("A" + $row + ":F" + $row)
    $destRow.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 2).Value = $subject
    $ws.Cells.Item($row, 4).Value = $logTime
    $ws.Cells.Item($row, 5).Value = $logType
    $ws.Cells.Item($row, 6).Value = $logUser

    $idCell = $ws.Cells.Item($row, 1)
    $idCell.NumberFormat = "@"
    $idCell.Value = $newRows[$row]

    $dateCell = $ws.Cells.Item($row, 3)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $logDate

    if ($row % 2 -eq 0) {
        $evenFormat.Copy() | Out-Null
        $idCell.PasteSpecial(-4122) | Out-Null
        $evenFormat.Copy() | Out-Null
        $dateCell.PasteSpecial(-4122) | Out-Null
    } else {
        $oddFormat.Copy() | Out-Null
        $idCell.PasteSpecial(-4122) | Out-Null
        $oddFormat.Copy() | Out-Null
        $dateCell.PasteSpecial(-4122) | Out-Null
    }
}
